# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# with refreshed values from the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.293.06"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.601.21"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'316.75"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'98.29"
$ws.Range("E6").Value = "  +4.09%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "'36.10"
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").Value = "'0.0817"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "'7.57"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "3.002.73"
$ws.Range("E13").Value = "  +3.52%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "2.655.28"
$ws.Range("E15").Value = "  +6.15%  "
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "'0.850"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "43.425.13"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").Value = "'12.78"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "0.0₃0971"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").Value = "'69.64"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").Value = "'255.13"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("E25").Value = "  +4.19%  "
$ws.Range("D26").Value = "'27.38"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'41.40"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "'10.32"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'156.44"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "'3.49"
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").Value = "'2.17"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("D35").Value = "'0.0812"
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").Value = "'18.82"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E39").Value = "  +10.78%  "
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("D41").Value = "'22.87"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("D42").Value = "'4.02"
$ws.Range("E42").Value = "  +6.93%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "2.015.57"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "'9.00"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "2.853.98"
$ws.Range("E48").Value = "  +3.72%  "
$ws.Range("D49").Value = "'83.67"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").Value = "'75.13"
$ws.Range("E50").Value = "  +2.51%  "
$ws.Range("D51").Value = "'0.195"
$ws.Range("E51").Value = "  +4.28%  "
